# Fruta / hortaliza, semanal
# Insert 2 new rows (this week's data) at the top of the data block (row 467),
# pushing the existing rows down by two. The newly inserted rows 467/468
# contain the new weekly price observations for Cebolla "1a (guarda)" /
# "2a (guarda)" at Vega Monumental Concepción.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 467:531 down to 469:533 by inserting two blank rows.
$ws.Rows("467:468").Insert()

# New row 467: Cebolla, Sin especificar, 1a (guarda)
$ws.Cells.Item(467,1).Value  = 11
$ws.Cells.Item(467,2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(467,3).Value  = "Bíobío"
$ws.Cells.Item(467,4).Value  = 44776
$ws.Cells.Item(467,5).Value  = 8
$ws.Cells.Item(467,6).Value  = 100112004
$ws.Cells.Item(467,7).Value  = "Cebolla"
$ws.Cells.Item(467,8).Value  = "Sin especificar"
$ws.Cells.Item(467,9).Value  = "1a (guarda)"
$ws.Cells.Item(467,10).Value = 200
$ws.Cells.Item(467,11).Value = 7000
$ws.Cells.Item(467,12).Value = 7000
$ws.Cells.Item(467,13).Value = 7000
$ws.Cells.Item(467,14).Value = "`$/malla 18 kilos"
$ws.Cells.Item(467,15).Value = "Región de O'Higgins"
$ws.Cells.Item(467,16).Value = 389
$ws.Cells.Item(467,17).Value = 18
$ws.Cells.Item(467,18).Value = "Hortaliza"

# New row 468: Cebolla, Sin especificar, 2a (guarda)
$ws.Cells.Item(468,1).Value  = 11
$ws.Cells.Item(468,2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(468,3).Value  = "Bíobío"
$ws.Cells.Item(468,4).Value  = 44776
$ws.Cells.Item(468,5).Value  = 8
$ws.Cells.Item(468,6).Value  = 100112004
$ws.Cells.Item(468,7).Value  = "Cebolla"
$ws.Cells.Item(468,8).Value  = "Sin especificar"
$ws.Cells.Item(468,9).Value  = "2a (guarda)"
$ws.Cells.Item(468,10).Value = 150
$ws.Cells.Item(468,11).Value = 6000
$ws.Cells.Item(468,12).Value = 6000
$ws.Cells.Item(468,13).Value = 6000
$ws.Cells.Item(468,14).Value = "`$/malla 18 kilos"
$ws.Cells.Item(468,15).Value = "Región de O'Higgins"
$ws.Cells.Item(468,16).Value = 333
$ws.Cells.Item(468,17).Value = 18
$ws.Cells.Item(468,18).Value = "Hortaliza"
